$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 from "E" to "P" (test status now Passed - "Completed with View Balance")
$ws.Range("B2").Value = "P"

# Add D2 with a literal empty-string value (not a blank cell). A plain
# Value = "" is indistinguishable from clearing the cell, so we go through
# Excel's text-prefix ( ' ) quoting to force a literal (empty) text value,
# then drop back to the Normal style so no stray quote-prefix formatting
# is left behind on the cell.
$ws.Range("D2").Value = "'"
$ws.Range("D2").Style = "Normal"
